$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-content permutation: the observation "records" spread across
# columns A,B,D-N,Q,R,Z,AB,AC were cyclically shuffled among certain
# rows (row numbers / other fixed columns such as S,T,U,V,W,Y,AD,AE,AG,
# AT,AW,AX,AY stay put). Below we set each destination row to the exact
# content it receives after the shuffle.

$rowData = @(
    @{ Row = 4; Cells = @{ "A" = 131009298; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589735; "R" = 6911227; "Z" = '10:12'; "AB" = '10:12'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 5; Cells = @{ "A" = 131009294; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589787; "R" = 6911183; "Z" = '10:24'; "AB" = '10:24'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 6; Cells = @{ "A" = 131009286; "B" = 80252; "D" = 'LC'; "E" = 6456; "F" = 'Skinnlav'; "G" = 'Leptogium saturninum'; "H" = '(Dicks.) Nyl.'; "I" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "Q" = 589865; "R" = 6911173; "Z" = '10:53'; "AB" = '10:53'; "AC" = $null } }
    @{ Row = 7; Cells = @{ "A" = 131009297; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589752; "R" = 6911214; "Z" = '10:18'; "AB" = '10:18'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 10; Cells = @{ "A" = 131009283; "B" = 79243; "D" = 'NT'; "E" = 6425; "F" = 'Garnlav'; "G" = 'Alectoria sarmentosa'; "H" = '(Ach.) Ach.'; "I" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "Q" = 589968; "R" = 6911120; "Z" = '11:03'; "AB" = '11:03'; "AC" = $null } }
    @{ Row = 16; Cells = @{ "A" = 131009301; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589700; "R" = 6911168; "Z" = '09:59'; "AB" = '09:59'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 17; Cells = @{ "A" = 131009282; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589986; "R" = 6911103; "Z" = '11:12'; "AB" = '11:12'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 21; Cells = @{ "A" = 131009498; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589629; "R" = 6911040; "Z" = '09:02'; "AB" = '09:02'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 23; Cells = @{ "A" = 131009304; "B" = 58043; "D" = 'NT'; "E" = 103021; "F" = 'Talltita'; "G" = 'Poecile montanus'; "H" = '(Conrad von Baldenstein, 1827)'; "I" = '1'; "K" = $null; "L" = $null; "M" = 'lockläte, övriga läten'; "N" = $null; "Q" = 589753; "R" = 6911167; "Z" = '09:45'; "AB" = '09:45'; "AC" = $null } }
    @{ Row = 24; Cells = @{ "A" = 131009305; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589757; "R" = 6911178; "Z" = '09:43'; "AB" = '09:43'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 25; Cells = @{ "A" = 131009295; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589741; "R" = 6911192; "Z" = '10:20'; "AB" = '10:20'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 26; Cells = @{ "A" = 131009269; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589671; "R" = 6911306; "Z" = '12:25'; "AB" = '12:25'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 27; Cells = @{ "A" = 131009271; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589680; "R" = 6911276; "Z" = '12:20'; "AB" = '12:20'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 28; Cells = @{ "A" = 131009281; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589899; "R" = 6911300; "Z" = '11:42'; "AB" = '11:42'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 29; Cells = @{ "A" = 131009285; "B" = 80252; "D" = 'LC'; "E" = 6456; "F" = 'Skinnlav'; "G" = 'Leptogium saturninum'; "H" = '(Dicks.) Nyl.'; "I" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "Q" = 589879; "R" = 6911153; "Z" = '10:56'; "AB" = '10:56'; "AC" = $null } }
    @{ Row = 30; Cells = @{ "A" = 131009302; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589716; "R" = 6911140; "Z" = '09:54'; "AB" = '09:54'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 31; Cells = @{ "A" = 131009306; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589791; "R" = 6911148; "Z" = '09:38'; "AB" = '09:38'; "AC" = 'färska och äldre ringhack på tall' } }
    @{ Row = 39; Cells = @{ "A" = 131009275; "B" = 57884; "D" = 'NT'; "E" = 100109; "F" = 'Tretåig hackspett'; "G" = 'Picoides tridactylus'; "H" = '(Linnaeus, 1758)'; "I" = $null; "K" = $null; "L" = $null; "M" = 'färska spår'; "N" = $null; "Q" = 589844; "R" = 6911365; "Z" = '11:53'; "AB" = '11:53'; "AC" = 'färska ringhack på tall' } }
    @{ Row = 40; Cells = @{ "A" = 131009291; "B" = 80252; "D" = 'LC'; "E" = 6456; "F" = 'Skinnlav'; "G" = 'Leptogium saturninum'; "H" = '(Dicks.) Nyl.'; "I" = $null; "K" = $null; "L" = $null; "M" = $null; "N" = $null; "Q" = 589791; "R" = 6911200; "Z" = '10:37'; "AB" = '10:37'; "AC" = $null } }
)

# Cells whose text content looks like a plain number (e.g. "1") need the
# target cell pre-formatted as Text, otherwise Excel auto-converts the
# assigned string into a numeric value.
$forceTextCells = @("I23")

foreach ($rd in $rowData) {
    $r = $rd.Row
    foreach ($col in $rd.Cells.Keys) {
        $addr = "$col$r"
        $val = $rd.Cells[$col]
        if ($forceTextCells -contains $addr) {
            $ws.Range($addr).NumberFormat = "@"
        }
        if ($val -eq $null) {
            $ws.Range($addr).Value = ""
        } else {
            $ws.Range($addr).Value = $val
        }
        if ($forceTextCells -contains $addr) {
            $ws.Range($addr).NumberFormat = "General"
        }
    }
}
